# Weekly reshuffle of "Hortaliza, Macroferia Regional de Talca - Arveja Verde":
# For each data row (2-33) the values in columns D,H,J,K,L,M,N,O,P,Q get
# replaced by the corresponding values taken from another row of the
# original sheet (a permutation of the existing rows - the other columns
# A,B,C,E,F,G,I,R are constant across every row so they are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by index) whose contents get shuffled between rows.
$cols = @(4, 8, 10, 11, 12, 13, 14, 15, 16, 17)   # D H J K L M N O P Q

# Mapping: destination row -> source row (values are copied from the
# source row's *original* contents into the destination row).
$rowMap = @{
    2 = 29
    3 = 17
    4 = 30
    5 = 28
    6 = 5
    7 = 9
    8 = 13
    9 = 4
    10 = 15
    11 = 19
    12 = 22
    13 = 12
    14 = 25
    15 = 2
    16 = 10
    17 = 20
    18 = 16
    19 = 11
    20 = 27
    21 = 24
    22 = 23
    23 = 7
    24 = 32
    25 = 6
    26 = 33
    27 = 3
    28 = 31
    29 = 18
    30 = 26
    31 = 14
    32 = 8
    33 = 21
}

# 1) Snapshot the original contents for every row/column involved before
#    writing anything, so later writes never clobber a value that is still
#    needed as a source for another row.
$snapshot = @{}
foreach ($r in $rowMap.Keys) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write the shuffled values back.
foreach ($r in $rowMap.Keys) {
    $srcRow = $rowMap[$r]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $srcVals[$c]
    }
}
